# Append new season rows (2024-25) of Flynn Dam depth readings, as per
# the "Added NC hikes (with updates to code in others)" commit.
#
# Existing data on Sheet1 runs through row 429 (last reading 2024-06-29).
# We extend it through row 436 with the next week's worth of readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date/time values use the workbook's existing "m/d/yy h:mm" (built-in
# numFmtId 22) style already applied to column B; depth values use the
# existing "0.00" (built-in numFmtId 2) style already applied to column C.
# Re-using those exact format strings keeps the new cells pointing at the
# same style records instead of minting new ones.
$dateFmt = "m/d/yy h:mm"
$numFmt = "0.00"

$rows = @(
    @{ Row = 430; Season = "2024-25"; DateTime = 45473.5625;          Depth = 2.07; Rain = 0.03 },
    @{ Row = 431; Season = "2024-25"; DateTime = 45474.693749999999;  Depth = 2.02 },
    @{ Row = 432; Season = "2024-25"; DateTime = 45475.651388888888;  Depth = 2;    Rain = 0.05 },
    @{ Row = 433; Season = "2024-25"; DateTime = 45476.453472222223;  Depth = 2;    Rain = 0.02 },
    @{ Row = 434; Season = "2024-25"; DateTime = 45477.530555555553;  Depth = 1.99 },
    @{ Row = 435; Season = "2024-25"; DateTime = 45478.490277777775;  Depth = 1.98; Rain = 0.01 },
    @{ Row = 436; Season = "2024-25"; DateTime = 45479.515972222223;  Depth = 1.96 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Season

    $ws.Cells.Item($row, 2).Value = $r.DateTime
    $ws.Cells.Item($row, 2).NumberFormat = $dateFmt

    $ws.Cells.Item($row, 3).Value = $r.Depth
    $ws.Cells.Item($row, 3).NumberFormat = $numFmt

    if ($r.ContainsKey("Rain")) {
        $ws.Cells.Item($row, 4).Value = $r.Rain
    }
}

# Mirror the author's final cursor position/selection on the bottom
# (frozen-below-row-3) pane after adding the new rows.
$ws.Range("D436").Select() | Out-Null
